$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Grid updates (A1:C5), driven by the "." markers ---
# Row 2: A2 clears (was "."), B2 becomes a single space, C2 becomes "."
$ws.Range("A2").ClearContents()
$ws.Range("B2").Value = " "
$ws.Range("C2").Value = "."

# Row 4: B4 becomes a single space (also gets a new style w/ explicit "no fill")
$ws.Range("B4").Value = " "
$ws.Range("B4").Interior.ColorIndex = -4142   # xlColorIndexNone -> forces applyFill on the xf

# --- Q column (static "paste values" snapshot of K column) ---
$ws.Range("Q4").Value = "{offsetx:0,offsety:3},"
$ws.Range("Q14").Value = ""

# --- Selection changes from Q1:Q15 to K1:K15 ---
$ws.Range("K1:K15").Select()
$excel.ActiveWindow.RangeSelection.Item(1).Activate()
